$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database values for cost_equity (X), roe_cost_equity (Y),
# cost_capital (AB), and roic_cost_capital (AC) columns on rows 2-4.

$ws.Range("X2").Value = 0.1066125655970851
$ws.Range("Y2").Value = -0.6244850329567142
$ws.Range("AB2").Value = 0.1043415746605175
$ws.Range("AC2").Value = -1.006867449551538

$ws.Range("X3").Value = 0.1031707526885612
$ws.Range("Y3").Value = 0.1424913477680597
$ws.Range("AB3").Value = 0.1029910568715701
$ws.Range("AC3").Value = -0.4617820006264095

$ws.Range("X4").Value = 0.1100543785056089
$ws.Range("Y4").Value = -1.391461413681488
$ws.Range("AB4").Value = 0.1056920924494648
$ws.Range("AC4").Value = -1.551952898476666
